$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I2").Value = 13
$ws.Range("J2").Value = 0.0009027777777777777
$ws.Range("K2").Value = 3690
$ws.Range("L2").Value = 0.00738
